$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New A (date serial) values for rows 2..53
$aVals = @(39400,39583,39765,39948,40130,40310,40494,40676,40862,41044,41228,41409,41592,41774,41957,42137,42321,42503,42689,42867,43053,43145,43235,43326,43418,43510,43600,43691,43783,43875,43966,44068,44159,44251,44341,44432,44525,44617,44706,44798,44890,44981,45071,45163,45254,45345,45436,45534,45618,45713,45800,45891)

# New B (data) values for rows 2..53
$bVals = @(-0.8957516004554691,1.910893826230975,-0.9623800849065276,-1.734840982136873,-1.4178810011895,-1.470147873721189,1.938408417700344,-0.163634158232469,0.4,-0.3,0,0.009235986179263023,0.3863474960573257,1.852186157158073,-0.1726927221574073,1.050807574684342,-0.2184803162966205,0.8513884674671885,0.4819278240608753,0.6629265129002277,1.164700738417963,0.3865811319001295,0.3857269132374052,0.4445645077595088,0.7,0.764905301728362,0.2,-0.1011524282971408,0.3399426662647187,-0.4401201376428787,-1.5,-7.231044133207007,4.665333711727925,-0.2672256391354182,-1.26956845314902,2.350117300589673,2.031122027629067,-0.5497577098212645,1.775201599534199,0.7199425862606432,0.5016525134024334,-0.5890802842202163,-1.005803602395417,0.5789791938159112,-0.3682387698544858,-0.2885446736570572,0,0.009541324901121584,1.028868575076984,1.007059757688495,-0.5319104727235242,0.4200898674779694)

for ($i = 0; $i -lt $aVals.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $aVals[$i]
    $ws.Cells.Item($r, 2).Value = $bVals[$i]
}

# Remove old rows 54..73 that no longer exist in the new data
$ws.Range("A54:B73").EntireRow.Delete() | Out-Null

